$d = $word.ActiveDocument

# Locate the start of the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph.
$startRange = $d.Content
$startFound = $startRange.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startRange.Expand(4)  # wdParagraph - expand to the whole paragraph

# Locate the end of the "... Creative Commons Attribution" paragraph, then also
# include the following (empty) paragraph so it disappears along with it.
$endRange = $d.Content
$endFound = $endRange.Find.Execute(
    "Creative Commons Attribution",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRange.Expand(4)     # wdParagraph - expand to the whole paragraph
$endRange.MoveEnd(1, 1) # wdCharacter - also swallow the following empty paragraph

if ($startFound -and $endFound) {
    $deleteRange = $d.Range($startRange.Start, $endRange.End)
    $deleteRange.Delete()
}
